# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns
# to the latest scraped snapshot, preserving each cell's original
# type/style. Price cells are forced to text (matching the source
# data, which stores prices like "1.00" or "64.152.10" as strings,
# not numbers) without leaving a residual explicit cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($range, $text) {
    # Force the cell to store $text verbatim as a string, even when it
    # looks numeric (e.g. "1.00"), then drop back to the default/unset
    # style so no explicit number-format style sticks to the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-PriceText $ws.Range("D2") '64.152.10'
$ws.Range("E2").Value = '  +5.66%  '
Set-PriceText $ws.Range("D3") '2.769.68'
$ws.Range("E3").Value = '  +4.87%  '
Set-PriceText $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  -0.02%  '
Set-PriceText $ws.Range("D5") '582.06'
$ws.Range("E5").Value = '  +1.09%  '
Set-PriceText $ws.Range("D6") '158.21'
$ws.Range("E6").Value = '  +10.01%  '
Set-PriceText $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  +0.16%  '
Set-PriceText $ws.Range("D8") '0.612'
$ws.Range("E8").Value = '  +2.17%  '
Set-PriceText $ws.Range("D9") '2.769.72'
$ws.Range("E9").Value = '  +4.21%  '
Set-PriceText $ws.Range("D10") '6.76'
$ws.Range("E10").Value = '  +2.96%  '
$ws.Range("E11").Value = '  +6.32%  '
Set-PriceText $ws.Range("D12") '0.395'
$ws.Range("E12").Value = '  +4.04%  '
$ws.Range("E13").Value = '  +1.98%  '
Set-PriceText $ws.Range("D14") '3.255.07'
$ws.Range("E14").Value = '  +4.66%  '
Set-PriceText $ws.Range("D15") '27.36'
$ws.Range("E15").Value = '  +4.68%  '
Set-PriceText $ws.Range("D16") '64.083.00'
$ws.Range("E16").Value = '  +5.59%  '
Set-PriceText $ws.Range("D17") '0.0000156'
$ws.Range("E17").Value = '  +9.01%  '
Set-PriceText $ws.Range("D18") '2.763.17'
$ws.Range("E18").Value = '  +3.98%  '
$ws.Range("E19").Value = '  +5.17%  '
Set-PriceText $ws.Range("D20") '4.97'
$ws.Range("E20").Value = '  +5.33%  '
Set-PriceText $ws.Range("D21") '363.42'
$ws.Range("E21").Value = '  +3.82%  '
Set-PriceText $ws.Range("D22") '7.03'
$ws.Range("E22").Value = '  +2.78%  '
Set-PriceText $ws.Range("D23") '1.00'
$ws.Range("E23").Value = '  +0.07%  '
Set-PriceText $ws.Range("D24") '0.537'
$ws.Range("E24").Value = '  +1.91%  '
Set-PriceText $ws.Range("D25") '67.27'
$ws.Range("E25").Value = '  +5.46%  '
Set-PriceText $ws.Range("D26") '0.172'
$ws.Range("E26").Value = '  +6.44%  '
Set-PriceText $ws.Range("D27") '8.63'
$ws.Range("E27").Value = '  +5.43%  '
Set-PriceText $ws.Range("D28") '0.997'
$ws.Range("E28").Value = '  -0.04%  '
Set-PriceText $ws.Range("D29") '0.0₃0914'
$ws.Range("E29").Value = '  +13.95%  '
Set-PriceText $ws.Range("D30") '2.04'
$ws.Range("E30").Value = '  +1.73%  '
Set-PriceText $ws.Range("D31") '7.19'
$ws.Range("E31").Value = '  +5.90%  '
Set-PriceText $ws.Range("D32") '1.27'
$ws.Range("E32").Value = '  +20.31%  '
Set-PriceText $ws.Range("D33") '175.30'
$ws.Range("E33").Value = '  +7.40%  '
Set-PriceText $ws.Range("D34") '0.999'
$ws.Range("E34").Value = '  +0.08%  '
Set-PriceText $ws.Range("D35") '20.74'
$ws.Range("E35").Value = '  +4.32%  '
$ws.Range("E36").Value = '  +6.54%  '
$ws.Range("E37").Value = '  +10.18%  '
$ws.Range("E38").Value = '  +9.86%  '
$ws.Range("E39").Value = '  +12.38%  '
Set-PriceText $ws.Range("D40") '4.33'
$ws.Range("E40").Value = '  +6.38%  '
Set-PriceText $ws.Range("D41") '340.69'
$ws.Range("E41").Value = '  +0.58%  '
Set-PriceText $ws.Range("D42") '39.33'
$ws.Range("E42").Value = '  +2.45%  '
Set-PriceText $ws.Range("D43") '5.84'
$ws.Range("E43").Value = '  +12.70%  '
Set-PriceText $ws.Range("D44") '22.05'
$ws.Range("E44").Value = '  +8.88%  '
Set-PriceText $ws.Range("D45") '22.31'
$ws.Range("E45").Value = '  +8.50%  '
Set-PriceText $ws.Range("D47") '0.651'
$ws.Range("E47").Value = '  +4.52%  '
Set-PriceText $ws.Range("D48") '0.0260'
$ws.Range("E48").Value = '  +4.31%  '
Set-PriceText $ws.Range("D49") '138.63'
$ws.Range("E49").Value = '  +4.40%  '
$ws.Range("E50").Value = '  +2.49%  '
Set-PriceText $ws.Range("D51") '1.00'
$ws.Range("E51").Value = '  +0.20%  '
